$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6666666666666666
$ws.Range("C2").Value = 0.7368421052631579
$ws.Range("D2").Value = 0.7
$ws.Range("B3").Value = 0.803921568627451
$ws.Range("C3").Value = 0.7454545454545455
$ws.Range("D3").Value = 0.7735849056603775
$ws.Range("B4").Value = 0.7419354838709677
$ws.Range("C4").Value = 0.7419354838709677
$ws.Range("D4").Value = 0.7419354838709677
$ws.Range("E4").Value = 0.7419354838709677
$ws.Range("B5").Value = 0.7352941176470589
$ws.Range("C5").Value = 0.7411483253588517
$ws.Range("D5").Value = 0.7367924528301888
$ws.Range("B6").Value = 0.7478389205144423
$ws.Range("C6").Value = 0.7419354838709677
$ws.Range("D6").Value = 0.7435179549604384
$ws.Range("B7").Value = 0.6666666666666666
$ws.Range("D7").Value = 0.7
$ws.Range("B8").Value = 0.803921568627451
$ws.Range("C8").Value = 0.7454545454545455
$ws.Range("D8").Value = 0.7735849056603775
$ws.Range("B9").Value = 0.7419354838709677
$ws.Range("C9").Value = 0.7419354838709677
$ws.Range("D9").Value = 0.7419354838709677
$ws.Range("E9").Value = 0.7419354838709677
$ws.Range("B10").Value = 0.7352941176470589
$ws.Range("C10").Value = 0.7411483253588517
$ws.Range("D10").Value = 0.7367924528301888
$ws.Range("B11").Value = 0.7478389205144423
$ws.Range("C11").Value = 0.7419354838709677
$ws.Range("D11").Value = 0.7435179549604384
$ws.Range("B12").Value = 0.65
$ws.Range("C12").Value = 0.3421052631578947
$ws.Range("D12").Value = 0.4482758620689655
$ws.Range("B13").Value = 0.6575342465753424
$ws.Range("C13").Value = 0.8727272727272727
$ws.Range("D13").Value = 0.7500000000000001
$ws.Range("B14").Value = 0.6559139784946236
$ws.Range("C14").Value = 0.6559139784946236
$ws.Range("D14").Value = 0.6559139784946236
$ws.Range("E14").Value = 0.6559139784946236
$ws.Range("B15").Value = 0.6537671232876712
$ws.Range("C15").Value = 0.6074162679425837
$ws.Range("D15").Value = 0.5991379310344829
$ws.Range("B16").Value = 0.6544557372219767
$ws.Range("C16").Value = 0.6559139784946236
$ws.Range("D16").Value = 0.6267148683722655
$ws.Range("B17").Value = 0.6129032258064516
$ws.Range("D17").Value = 0.5507246376811595
$ws.Range("B18").Value = 0.6935483870967742
$ws.Range("C18").Value = 0.7818181818181819
$ws.Range("D18").Value = 0.7350427350427351
$ws.Range("B19").Value = 0.6666666666666666
$ws.Range("C19").Value = 0.6666666666666666
$ws.Range("D19").Value = 0.6666666666666666
$ws.Range("E19").Value = 0.6666666666666666
$ws.Range("B20").Value = 0.653225806451613
$ws.Range("C20").Value = 0.6409090909090909
$ws.Range("D20").Value = 0.6428836863619474
$ws.Range("B21").Value = 0.660596600763094
$ws.Range("C21").Value = 0.6666666666666666
$ws.Range("D21").Value = 0.6597299640777903
$ws.Range("B22").Value = 0.6458333333333334
$ws.Range("C22").Value = 0.8157894736842105
$ws.Range("D22").Value = 0.7209302325581395
$ws.Range("B23").Value = 0.8444444444444444
$ws.Range("C23").Value = 0.6909090909090909
$ws.Range("D23").Value = 0.7599999999999999
$ws.Range("B24").Value = 0.7419354838709677
$ws.Range("C24").Value = 0.7419354838709677
$ws.Range("D24").Value = 0.7419354838709677
$ws.Range("E24").Value = 0.7419354838709677
$ws.Range("B25").Value = 0.7451388888888889
$ws.Range("C25").Value = 0.7533492822966508
$ws.Range("D25").Value = 0.7404651162790696
$ws.Range("B26").Value = 0.7632915173237754
$ws.Range("C26").Value = 0.7419354838709677
$ws.Range("D26").Value = 0.7440360090022505
